# Updated test data as per new implementation
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# Rename the column header labels (shared strings updated in place)
$ws.Range("L8").Value = "Alarm Current(A)"
$ws.Range("M8").Value = "Standby Current(A)"

# Update the saved selection/active cell on the "Add Panels" sheet
$ws.Activate()
$ws.Range("L8:M8").Select()
